$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Remove the "Money" entity type (row 11 on the "Entity types" sheet) and
# shift the rows below it up.
$ws1.Rows.Item(11).Delete()

# Re-alphabetize the remaining 24 entity types in column A.
$sorted = @(
  "Anatomy",
  "Award",
  "Broadcaster",
  "Company",
  "Crime",
  "Drug",
  "EmailAddress",
  "Facility",
  "GeographicFeature",
  "HealthCondition",
  "Hashtag",
  "IPAddress",
  "JobTitle",
  "Location",
  "Movie",
  "MusicGroup",
  "NaturalEvent",
  "Organization",
  "Person",
  "Sport",
  "SportingEvent",
  "TelevisionShow",
  "TwitterHandle",
  "Vehicle"
)
for ($i = 0; $i -lt $sorted.Length; $i++) {
  $ws1.Cells.Item($i + 1, 1).Value = $sorted[$i]
}

# Darken the font color of the list to a dark gray, matching the new look.
$ws1.Range("A1:A24").Font.Color = 3355443

# Resize column A to fit the newly sorted (and slightly different) content.
$ws1.Columns.Item(1).AutoFit()

# Switch the page to portrait orientation on both sheets.
$ws1.PageSetup.Orientation = 1
$ws2.PageSetup.Orientation = 1

# Restore the on-screen selections recorded for each sheet.
$ws2.Activate()
$ws2.Range("A14").Select()

$ws1.Activate()
$ws1.Range("C8").Select()
